# Auto update Excel log
# Append newly recorded sensor readings (2026-02-01, ~20:40-20:41) to the
# mmWave(InBed), mmWave(BR) and mmWave(HR) worksheets, matching the live feed export.

$wb = $excel.ActiveWorkbook

# --- mmWave(InBed) : append rows 30-41 ---
$ws = $wb.Worksheets.Item("mmWave(InBed)")
$newRows = @(
    @("2026-02-01", "20:40:47", "20:00", "Bedroom", "Out of Bed", "Empty"),
    @("2026-02-01", "20:40:48", "20:00", "Bedroom", "In Bed", "Occupied"),
    @("2026-02-01", "20:40:48", "20:00", "Bedroom", "In Bed", "Occupied"),
    @("2026-02-01", "20:40:51", "20:00", "Bedroom", "In Bed", "Occupied"),
    @("2026-02-01", "20:40:51", "20:00", "Bedroom", "In Bed", "Occupied"),
    @("2026-02-01", "20:40:53", "20:00", "Bedroom", "In Bed", "Occupied"),
    @("2026-02-01", "20:40:53", "20:00", "Bedroom", "In Bed", "Occupied"),
    @("2026-02-01", "20:40:54", "20:00", "Bedroom", "In Bed", "Occupied"),
    @("2026-02-01", "20:41:07", "20:00", "Bedroom", "In Bed", "Occupied"),
    @("2026-02-01", "20:41:08", "20:00", "Bedroom", "In Bed", "Occupied"),
    @("2026-02-01", "20:41:08", "20:00", "Bedroom", "In Bed", "Occupied"),
    @("2026-02-01", "20:41:13", "20:00", "Bedroom", "In Bed", "Occupied")
)
$r = 30
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 1).ClearFormats()
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 2).ClearFormats()
    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 3).ClearFormats()
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# --- mmWave(BR) : append rows 30-39 ---
$ws = $wb.Worksheets.Item("mmWave(BR)")
$newRows = @(
    @("2026-02-01", "20:40:49", "20:00", "Bedroom", 2, "Occupied"),
    @("2026-02-01", "20:40:51", "20:00", "Bedroom", 18, "Occupied"),
    @("2026-02-01", "20:40:52", "20:00", "Bedroom", 4, "Occupied"),
    @("2026-02-01", "20:40:53", "20:00", "Bedroom", 2, "Occupied"),
    @("2026-02-01", "20:40:54", "20:00", "Bedroom", 29, "Occupied"),
    @("2026-02-01", "20:40:55", "20:00", "Bedroom", 2, "Occupied"),
    @("2026-02-01", "20:41:07", "20:00", "Bedroom", 16, "Occupied"),
    @("2026-02-01", "20:41:08", "20:00", "Bedroom", 5, "Occupied"),
    @("2026-02-01", "20:41:09", "20:00", "Bedroom", 2, "Occupied"),
    @("2026-02-01", "20:41:13", "20:00", "Bedroom", 1, "Occupied")
)
$r = 30
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 1).ClearFormats()
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 2).ClearFormats()
    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 3).ClearFormats()
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# --- mmWave(HR) : append rows 30-39 ---
$ws = $wb.Worksheets.Item("mmWave(HR)")
$newRows = @(
    @("2026-02-01", "20:40:49", "20:00", "Bedroom", 50, "Occupied"),
    @("2026-02-01", "20:40:51", "20:00", "Bedroom", 66, "Occupied"),
    @("2026-02-01", "20:40:52", "20:00", "Bedroom", 52, "Occupied"),
    @("2026-02-01", "20:40:53", "20:00", "Bedroom", 50, "Occupied"),
    @("2026-02-01", "20:40:53", "20:00", "Bedroom", 77, "Occupied"),
    @("2026-02-01", "20:40:55", "20:00", "Bedroom", 50, "Occupied"),
    @("2026-02-01", "20:41:07", "20:00", "Bedroom", 64, "Occupied"),
    @("2026-02-01", "20:41:08", "20:00", "Bedroom", 53, "Occupied"),
    @("2026-02-01", "20:41:09", "20:00", "Bedroom", 50, "Occupied"),
    @("2026-02-01", "20:41:13", "20:00", "Bedroom", 49, "Occupied")
)
$r = 30
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 1).ClearFormats()
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 2).ClearFormats()
    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 3).ClearFormats()
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}
